$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'256.22"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'-0.49%"
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'26.86"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'-0.59%"
$ws.Range("E3").ClearFormats()
$ws.Range("D4").Value = "'4.506"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'-5.02%"
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'0.05881"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'-1.45%"
$ws.Range("E5").ClearFormats()
$ws.Range("E6").Value = "'-1.06%"
$ws.Range("E6").ClearFormats()
$ws.Range("D7").Value = "'0.8502"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'-2.32%"
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'0.9301"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'-1.73%"
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'0.1377"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'-2.10%"
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'0.04526"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'24.85%"
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'0.07029"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'-2.22%"
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'0.03063"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'-3.48%"
$ws.Range("E12").ClearFormats()
$ws.Range("E13").Value = "'-1.65%"
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'0.001527"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'-1.53%"
$ws.Range("E14").ClearFormats()
$ws.Range("B15").Value = "'One"
$ws.Range("B15").ClearFormats()
$ws.Range("C15").Value = "'https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("C15").ClearFormats()
$ws.Range("D15").Value = "'0.0006028"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'-1.56%"
$ws.Range("E15").ClearFormats()
$ws.Range("B16").Value = "'TigerCash"
$ws.Range("B16").ClearFormats()
$ws.Range("C16").Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("C16").ClearFormats()
$ws.Range("D16").Value = "'0.006010"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'0.26%"
$ws.Range("E16").ClearFormats()
$ws.Range("B17").Value = "'LEO"
$ws.Range("B17").ClearFormats()
$ws.Range("C17").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("C17").ClearFormats()
$ws.Range("D17").Value = "'3.482"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'-0.48%"
$ws.Range("E17").ClearFormats()
$ws.Range("B18").Value = "'GateToken"
$ws.Range("B18").ClearFormats()
$ws.Range("C18").Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("C18").ClearFormats()
$ws.Range("D18").Value = "'3.171"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'-0.18%"
$ws.Range("E18").ClearFormats()
$ws.Range("B19").Value = "'BTSEToken"
$ws.Range("B19").ClearFormats()
$ws.Range("C19").Value = "'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("C19").ClearFormats()
$ws.Range("D19").Value = "'2.204"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'-1.61%"
$ws.Range("E19").ClearFormats()
$ws.Range("E21").Value = "'-1.63%"
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'3.922"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'2.80%"
$ws.Range("E22").ClearFormats()
$ws.Range("E23").Value = "'0.74%"
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'0.001221"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'-0.58%"
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'0.004304"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'-4.35%"
$ws.Range("E25").ClearFormats()
$ws.Range("E26").Value = "'-0.04%"
$ws.Range("E26").ClearFormats()
$ws.Range("E27").Value = "'1.94%"
$ws.Range("E27").ClearFormats()
$ws.Range("E40").Value = "'-0.51%"
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = "'0.006288"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'1.41%"
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'0.1098"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'-0.34%"
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'0.002199"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'-2.36%"
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'0.01384"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'30.53%"
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'0.00005335"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'-3.07%"
$ws.Range("E45").ClearFormats()
$ws.Range("D48").Value = "'0.2519"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'10,970.15%"
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("D49").ClearFormats()
$ws.Range("D50").Value = "'0.0001999"
$ws.Range("D50").ClearFormats()
